$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N, shifting N->O, O->P, P->Q
$ws3.Columns("N:N").Insert() | Out-Null
$ws3.Columns("N:N").ColumnWidth = 10.17

# Switch active sheet to "Repayment schedule" and select S7
$ws3.Activate() | Out-Null
$ws3.Range("S7").Select() | Out-Null
